# Added test for adding a book to user profile
#
# 1) LoginCredentials!A2:B2 gets a new valid username/password pair
#    (B2 also carries a hyperlink, matching the pre-existing Hyperlink
#    cell style that already sat, empty, on B2).
# 2) A new "Books" worksheet is added after LoginCredentials with a
#    single book title, styled in a custom monospace font.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- LoginCredentials: new credential row -------------------------------
$ws1.Range("A2").Value = "dragoljubqa"
$ws1.Range("B2").Value = "Qwerty123!@#"

$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://www.packtpub.com/product/speaking-javascript/9781449365029")
$ws1.Range("B2").Style = "Hyperlink"

$ws1.Range("B2").Select()

# --- New "Books" worksheet ------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Books"

$ws2.Range("A1").Value = "Speaking JavaScript"

$bookFont = $ws2.Range("A1").Font
$bookFont.Name = "JetBrains Mono"
$bookFont.Family = 3
$bookFont.Size = 10
$bookFont.Color = 0x73AB6A
$ws2.Range("A1").VerticalAlignment = -4108

$ws2.Columns.Item(1).AutoFit()
